$d = $word.ActiveDocument

# --- Insert the new paragraph "Vilken fin värld hej världen" right after
#     paragraph 1 (before the two blank paragraphs that precede "TIIIM"). ---
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(2)

# Write the text, plus one throw-away trailing character. We need that
# extra character because this COM host mis-resolves a *collapsed*
# (zero-length) Range positioned exactly at "end of paragraph text, just
# before the paragraph mark" -> Bookmarks.Add anchored there lands back at
# document position 0 instead of where it should be. Writing one extra
# character first turns the bookmark anchor point into an ordinary
# mid-text position (which resolves correctly), and we delete that
# character again afterwards.
$placeholder = "x"
$p2.Range.Text = "Vilken fin värld hej världen" + $placeholder

# Collapsed anchor sitting right before the placeholder character (i.e.
# exactly where the bookmark should end up: right after the real text).
$anchorPos = $p2.Range.End - (1 + $placeholder.Length)
$bmRange = $d.Range($anchorPos, $anchorPos)

# Re-adding a bookmark named "_GoBack" removes any previous bookmark of
# the same name first (Word only allows one bookmark per name), which
# takes care of stripping the bookmarkStart/bookmarkEnd pair that used to
# sit in the "TIIIM" paragraph, exactly as the diff wants.
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the placeholder character now that the bookmark is anchored.
$placeholderRange = $d.Range($p2.Range.End - (1 + $placeholder.Length), $p2.Range.End - 1)
$placeholderRange.Delete()

Write-Output "done"
